# Locate the paragraph that currently ends with "...programming lang"
# (3rd paragraph in git.docx) and replace its contents with the
# post-diff OOXML: the run "lang" is split into its own run wrapped in
# spellStart/gramStart..spellEnd/gramEnd proofErr marks, and a new
# paragraph "Git is very useful" is appended after it (carrying forward
# the trailing bookmarkStart/bookmarkEnd "_GoBack" pair).

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*can be used with any other programming lang*") {
        $target = $p
    }
}

$r = $target.Range
# Trim the trailing paragraph mark from the range so InsertXML only
# replaces the paragraph's run content (the mark's own formatting /
# the bookmark living right before it is re-created explicitly below).
$r.End = $r.End - 1

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00644499" w:rsidRDefault="00644499"><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> can be used with any other programming </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>lang</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is very useful</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)
